$d = $word.ActiveDocument

$find = " In case of error, the timestamp is also sent."
$replace = ""

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replace, 2)
